$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New rows of service records (DATE, VEHICLE REG NO, VEHICLE BRAND, ISSUE, STATUS, AMOUNT, CASH TYPE)
$rows = @(
    @{A=44818; B="MH02AP8026"; C="IKON";        D="PMS";                E="WORK DONE DELIVERED"; F=4427;  G="CARD"},
    @{A=44818; B="KA03NC7377"; C="ECOSPORT";     D="GENERAL CHECKUP";    E="WORK DONE "},
    @{A=44818; B="KA53MC3710"; C="SWFIT DZIRE";  D="PMS";                E="WORK DONE DELIVERED"; F=6209;  G="GPAY"},
    @{A=44818; B="HR26CC3104"; C="LINEA";        D="RUNNING REPAIR";     E="WORK IN PROGRESS"},
    @{A=44818; B="KA53M1408";  C="ACCENT";       D="PMS";                E="WORK DONE DELIVERED"; F=5300;  G="PAYTM"},
    @{A=44818; B="KA03MQ3773 ";C="VENTO";        D="BODY SHOP";          E="WORK DONE DELIVERED"; F=30000; G="CARD"},
    @{A=44818; B="TN64L3591";  C="DATSUN";       D="CLUTCH PROBLEM";     E="WORK DONE DELIVERED"; F=12750; G="G PAY"},
    @{A=44818; B="KA05MP4915"; C="PUNTO";        D="DRIVE SHAFT CHANGE"; E="WORK DONE DELIVERED"; F=14538; G="CREDIT"},
    @{A=44818; B="KA53MG5551"; C="BENZ";         D="BODY SHOP";          E="WORK DONE DELIVERED"; F=11000}
)

$startRow = 401
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Range("A$r").NumberFormat = "[$-409]d/mmm/yyyy;@"

    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E

    if ($row.ContainsKey("F")) {
        $ws.Cells.Item($r, 6).Value = $row.F
        $ws.Range("F$r").NumberFormat = '"₹"#,##0;"₹"\-#,##0'
    }
    if ($row.ContainsKey("G")) {
        $ws.Cells.Item($r, 7).Value = $row.G
    }
}

$excel.ActiveWindow.ScrollRow = 385
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G409").Select()
